# #73 Update the normalization
# Applies the changes described by the diff:
#  - Products sheet: B13 formula gets a new "Part_ID" column reference in the
#    generated SQL text, and rows 14:19 (the rest of the old shared-formula
#    fill) are removed.
#  - ProductFamily sheet: selection collapses from B9:B13 to B9 (cosmetic).
#  - Orders sheet: a new column H is added with a CONCAT formula that builds
#    an "insert into Customer_Order_Product" statement for every order row.
#  - Workstations becomes the active/selected sheet (was BOM before).
#  - BOM sheet is scrolled down and is no longer the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Products sheet
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")

# Update the remaining formula (row 13) to include the Part_ID column name.
$wsProducts.Range("B13").Formula = '=CONCAT("insert into Product (Product_ID, Product_Name, Product_Description, Factory_Plant_ID, Family_ID, Part_ID) values (''",A2,"'', ''",B2,"'', ''",C2,"'', 1, ",D2,");")'

# Remove the old rows 14-19 (they duplicated the formula for other products).
$wsProducts.Rows("14:19").Delete()

# Restore the selection to match the new, smaller used range.
$wsProducts.Range("B14:B19").Select()

# ---------------------------------------------------------------------
# 2. ProductFamily sheet
# ---------------------------------------------------------------------
$wsProductFamily = $wb.Worksheets.Item("ProductFamily")
$wsProductFamily.Range("B9").Select()

# ---------------------------------------------------------------------
# 3. Orders sheet
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")

$wsOrders.Range("H2").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A2,", ''",C2,"'', ",D2,");")'
$wsOrders.Range("H3").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A3,", ''",C3,"'', ",D3,");")'
$wsOrders.Range("H4").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A4,", ''",C4,"'', ",D4,");")'
$wsOrders.Range("H5").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A5,", ''",C5,"'', ",D5,");")'
$wsOrders.Range("H6").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A6,", ''",C6,"'', ",D6,");")'
$wsOrders.Range("H7").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A7,", ''",C7,"'', ",D7,");")'
$wsOrders.Range("H8").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A8,", ''",C8,"'', ",D8,");")'
$wsOrders.Range("H9").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A9,", ''",C9,"'', ",D9,");")'
$wsOrders.Range("H10").Formula = '=CONCAT("insert into Customer_Order_Product(Customer_Order_ID, Product_ID, Quantity) values (",A10,", ''",C10,"'', ",D10,");")'

$wsOrders.Range("I17").Select()

# ---------------------------------------------------------------------
# 4. BOM sheet - scroll down, no longer the active tab
# ---------------------------------------------------------------------
$wsBOM = $wb.Worksheets.Item("BOM")
$wsBOM.Activate()
$wsBOM.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 6
$wsBOM.Range("B60").Select()

# ---------------------------------------------------------------------
# 5. Workstations sheet becomes the active/selected sheet
# ---------------------------------------------------------------------
$wsWorkstations = $wb.Worksheets.Item("Workstations")
$wsWorkstations.Activate()
$wsWorkstations.Range("A23").Select()
